# Apply the "Add removing element from slice algorithm explanation" edit:
#  - Add a gopher picture to slides 3, 4, 5 (bottom-left corner)
#  - Fix slide 6 title run split, add full body text + picture
#  - Add a brand-new slide 7 ("Slice Operator & Removing Elements") with body text + picture

$p = $ppt.ActivePresentation

function Add-GopherPicture {
    param(
        $Slide,
        [double]$X,
        [double]$Y
    )
    # Copy the existing gopher picture (already embedded on slide 1) so the
    # media relationship / blip is reused instead of creating a brand-new part.
    $srcSlide = $p.Slides.Item(1)
    $srcPic = $srcSlide.Shapes.Item(3)
    $srcPic.Copy()
    $newPic = $Slide.Shapes.Paste()
    $newPic.Name = "Content Placeholder 4"
    $newPic.Left = $X / 12700.0
    $newPic.Top = $Y / 12700.0
    return $newPic
}

# ---------------------------------------------------------------------------
# Slide 3 ("Side note on Go arrays") - add gopher picture
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
Add-GopherPicture -Slide $slide3 -X 1007532 -Y 3824685 | Out-Null

# ---------------------------------------------------------------------------
# Slide 4 ("Limitations on Arrays") - add gopher picture
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
Add-GopherPicture -Slide $slide4 -X 1150777 -Y 3781142 | Out-Null

# ---------------------------------------------------------------------------
# Slide 5 ("Slice internals") - add gopher picture
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
Add-GopherPicture -Slide $slide5 -X 1107234 -Y 3698748 | Out-Null

# ---------------------------------------------------------------------------
# Slide 6 ("Slice Syntax & Functions") - fix title run, fill body, add picture
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

# Merge the title's two runs ("Slice Syntax " + "& Functions") into a single run.
$title6 = $slide6.Shapes.Item(1)
$title6.TextFrame.TextRange.Text = "placeholder"
$title6.TextFrame.TextRange.Text = "Slice Syntax & Functions"

# Fill in the body placeholder with the full explanation text.
$body6 = $slide6.Shapes.Item(2)
$body6.TextFrame.TextRange.Text = "Declare a slice variable (Two Ways)`r" + `
    "variableName := make([]type, length, *capacity)`r" + `
    "variableName := []type{initial elements}`r" + `
    "Add new element to slice`r" + `
    "Use built in function append - newSlice = append(oldSlice, newElement). Remember append function always returns a new slice.`r" + `
    "Get length of slice`r" + `
    "Use built in function len - len(mySlice).`r" + `
    "Copy one slice variable to another slice variable`r" + `
    "To copy one slice to another use the build in function copy - copy(dest []type, src []type).`r" + `
    " "

$tr6 = $body6.TextFrame.TextRange

# Paragraph 10 (last, blank) -> indent level 2 (lvl="1"), empty text already there.
$para10 = $tr6.Paragraphs(10, 1)
$para10.IndentLevel = 2
$para10.Text = ""

# Paragraph 9: "To copy one slice ... copy(" + "dest" + " []type, " + "src" + " []type)."
$para9 = $tr6.Paragraphs(9, 1)
$para9.IndentLevel = 2
$para9.Text = "To copy one slice to another use the build in function copy – copy("
$para9.InsertAfter("dest")
$para9.InsertAfter(" []type, ")
$para9.InsertAfter("src")
$para9.InsertAfter(" []type).")

# Paragraph 8: "Copy one slice variable to another slice variable"
$para8 = $tr6.Paragraphs(8, 1)
$para8.IndentLevel = 1

# Paragraph 7: "Use built in function len - len(mySlice)."
$para7 = $tr6.Paragraphs(7, 1)
$para7.IndentLevel = 2
$para7.Text = "Use built in function "
$para7.InsertAfter("len")
$para7.InsertAfter(" – ")
$para7.InsertAfter("len")
$para7.InsertAfter("(")
$para7.InsertAfter("mySlice")
$para7.InsertAfter(").")

# Paragraph 6: "Get length of slice"
$para6 = $tr6.Paragraphs(6, 1)
$para6.IndentLevel = 1

# Paragraph 5: append(...) explanation
$para5 = $tr6.Paragraphs(5, 1)
$para5.IndentLevel = 2
$para5.Text = "Use built in function append – "
$para5.InsertAfter("newSlice")
$para5.InsertAfter(" = append(")
$para5.InsertAfter("oldSlice")
$para5.InsertAfter(", ")
$para5.InsertAfter("newElement")
$para5.InsertAfter("). Remember append function always returns a new slice.")

# Paragraph 4: "Add new element to slice"
$para4 = $tr6.Paragraphs(4, 1)
$para4.IndentLevel = 1

# Paragraph 3: "variableName := []type{initial elements}"
$para3 = $tr6.Paragraphs(3, 1)
$para3.IndentLevel = 2
$para3.Text = "variableName"
$para3.InsertAfter(" := []type{initial elements}")

# Paragraph 2: "variableName := make([]type, length, *capacity)"
$para2 = $tr6.Paragraphs(2, 1)
$para2.IndentLevel = 2
$para2.Text = "variableName"
$para2.InsertAfter(" := make([]type, length, *capacity)")

# Paragraph 1 stays at default indent level (0) - "Declare a slice variable (Two Ways)"

Add-GopherPicture -Slide $slide6 -X 1119675 -Y 3793583 | Out-Null

# ---------------------------------------------------------------------------
# Slide 7 (new) - "Slice Operator & Removing Elements"
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Add(7, 2)

$title7 = $slide7.Shapes.Item(1)
$title7.TextFrame.TextRange.Text = "Slice Operator & Removing Elements"

$body7 = $slide7.Shapes.Item(2)
$body7.TextFrame.TextRange.Text = "variableName[fromLeft:toRight]`r" + `
    "For example we have the following slice [0, 1, 2, 3, 4, 5]`r" + `
    "If we write the following slice[0:4], we get returned a new slice containing the elements from index 0 (inclusive) to index 4 (not inclusive). So in our example we get a slice containing [0, 1, 2, 3].`r" + `
    "The following algorithm can be used to remove an element from a slice at a certain index.`r" + `
    "Slice = append(Slice[:indexOfElementToRemove], Slice[indexOfElementToRemove + 1:]...)"

$tr7 = $body7.TextFrame.TextRange

# Paragraph 5: Slice = append(Slice[:...], Slice[... + 1:]...)
$para7_5 = $tr7.Paragraphs(5, 1)
$para7_5.Text = "Slice = append(Slice[:"
$para7_5.InsertAfter("indexOfElementToRemove")
$para7_5.InsertAfter("], Slice[")
$para7_5.InsertAfter("indexOfElementToRemove")
$para7_5.InsertAfter(" + 1:]…)")

# Paragraph 1: variableName[fromLeft:toRight]
$para7_1 = $tr7.Paragraphs(1, 1)
$para7_1.Text = "variableName"
$para7_1.InsertAfter("[")
$para7_1.InsertAfter("fromLeft:toRight")
$para7_1.InsertAfter("]")

Add-GopherPicture -Slide $slide7 -X 1200540 -Y 4135706 | Out-Null
